$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.050.78"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.792.48"
$ws.Range("E3").Value = "  -2.33%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.007"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.80"
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4177"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3563"
$ws.Range("E8").Value = "  -2.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07054"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8429"
$ws.Range("E10").Value = "  -2.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.08"
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.853.37"
$ws.Range("E12").Value = "  -8.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.256"
$ws.Range("E13").Value = "  -2.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.338"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06861"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.82"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008709"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.02"
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.305.18"
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.047"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.65"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.052.67"
$ws.Range("E24").Value = "  -12.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.961"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.05"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.14"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.002"
$ws.Range("E28").Value = "  -4.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.42"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.659"
$ws.Range("E30").Value = "  -9.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08886"
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.886"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.341"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.008"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.080"
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.073"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01891"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05104"
$ws.Range("E39").Value = "  -4.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4957"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.667"
$ws.Range("E42").Value = "  -5.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.150"
$ws.Range("E43").Value = "  -9.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.996"
$ws.Range("E44").Value = "  -5.71%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.66"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.007"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.19"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06309"
$ws.Range("E48").Value = "  -3.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4525"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.590"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.10"
$ws.Range("E51").Value = "  -3.47%  "
